$wb = $excel.ActiveWorkbook

# --- Sheet 1: "All Orders" ---------------------------------------------
$ws1 = $wb.Worksheets.Item("All Orders")

# Insert a new row above the existing data row (old row 2 shifts to row 3)
$ws1.Range("A2:N2").Insert()

# Populate the new order in row 2
$ws1.Range("A2").Value2 = 2
$ws1.Range("B2").Value2 = "2026-01-13 10:20"
$ws1.Range("C2").Value2 = "Pooja"
$ws1.Range("D2").Value2 = "A 1608"
$ws1.Range("E2").Value2 = ""
$ws1.Range("F2").Value2 = "Wheat Chapati x1"
$ws1.Range("G2").Value2 = 15
$ws1.Range("H2").Value2 = "NEW"
$ws1.Range("I2").Value2 = "PENDING"
$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value2 = "2026-01-13"
$ws1.Range("K2").Value2 = "15:50"
$ws1.Range("L2").Value2 = ""
$ws1.Range("M2").Value2 = ""
$ws1.Range("N2").Value2 = ""

# --- Sheet 2: "Daily Summary" -------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")

$ws2.Range("B2").Value2 = 2
$ws2.Range("E2").Value2 = 45
$ws2.Range("G2").Value2 = 45
